# Update the "Fanano" report with data through 9 August 2021 (rows 329-343).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date (column A), nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$data = @(
    @(329, 44403, 0, 0, 0),
    @(330, 44404, 0, 0, 0),
    @(331, 44405, 0, 0, 0),
    @(332, 44406, 0, 0, 0),
    @(333, 44407, 0, 0, 0),
    @(334, 44408, 0, 0, 0),
    @(335, 44409, 0, 0, 0),
    @(336, 44410, 1, 1, 33.71544167228591),
    @(337, 44411, 1, 2, 67.43088334457181),
    @(338, 44412, 1, 3, 101.1463250168577),
    @(339, 44413, 1, 4, 134.8617666891436),
    @(340, 44414, 2, 6, 202.2926500337155),
    @(341, 44415, 4, 10, 337.1544167228591),
    @(342, 44416, 7, 17, 573.1625084288604),
    @(343, 44417, 1, 17, 573.1625084288604)
)

# Row 328 carries the "date" style (border + center + yyyy-mm-dd hh:mm:ss
# number format) that every column-A cell in this sheet uses; clone it down
# onto the new rows' A cells before writing values.
$styleSource = $ws.Range("A328")

foreach ($r in $data) {
    $row = $r[0]
    $styleSource.Copy()
    $ws.Range("A$row").PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

$excel.CutCopyMode = 0
